$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Column width adjustments (values pre-compensated for the ~5/6
#    character offset Excel adds when converting ColumnWidth -> stored
#    "width" so the saved OOXML <col width="..."> lands on the target
#    integer exactly).
# -----------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth  = 30.166666666666668   # 37 -> 31
$ws.Columns.Item(3).ColumnWidth  = 31.166666666666668   # 39 -> 32
$ws.Columns.Item(4).ColumnWidth  = 30.166666666666668   # 39 -> 31
$ws.Columns.Item(5).ColumnWidth  = 29.166666666666668   # 39 -> 30
$ws.Columns.Item(6).ColumnWidth  = 30.166666666666668   # 39 -> 31
$ws.Columns.Item(7).ColumnWidth  = 31.166666666666668   # 39 -> 32
$ws.Columns.Item(16).ColumnWidth = 47.166666666666664   # 40 -> 48
$ws.Columns.Item(19).ColumnWidth = 43.166666666666664   # 37 -> 44
$ws.Columns.Item(20).ColumnWidth = 33.166666666666664   # 44 -> 34
$ws.Columns.Item(21).ColumnWidth = 30.166666666666668   # 34 -> 31
$ws.Columns.Item(22).ColumnWidth = 59.166666666666664   # 36 -> 60
$ws.Columns.Item(23).ColumnWidth = 36.166666666666664   # 31 -> 37
$ws.Columns.Item(24).ColumnWidth = 33.166666666666664   # 33 -> 34
$ws.Columns.Item(25).ColumnWidth = 30.166666666666668   # 60 -> 31

# -----------------------------------------------------------------------
# 2) Header row (row 1) relabeling - rename the Katalon-AI field names.
# -----------------------------------------------------------------------
$ws.Range("A1").Value = "div_backdropElements_class"
$ws.Range("B1").Value = "div_testCases_internalText"
$ws.Range("C1").Value = "div_testCases_internalText_1"
$ws.Range("D1").Value = "div_testCases_internalText_2"
$ws.Range("E1").Value = "div_testCases_internalText_3"
$ws.Range("F1").Value = "div_testCases_internalText_4"
$ws.Range("G1").Value = "div_testCases_internalText_5"
$ws.Range("I1").Value = "div_testSuites_internalText"
$ws.Range("J1").Value = "div_testSuites_internalText_1"
$ws.Range("K1").Value = "input_rowSelectionCheckbox_class"
$ws.Range("L1").Value = "input_rowSelectionCheckbox_class_1"
$ws.Range("M1").Value = "input_rowSelectionCheckbox_class_2"
$ws.Range("N1").Value = "input_rowSelectionCheckbox_class_3"
$ws.Range("O1").Value = "input_rowSelectionCheckbox_class_4"
$ws.Range("P1").Value = "input_rowSelectionCheckbox_internalRoleRowName"
$ws.Range("Q1").Value = "input_rowSelectionCheckbox_internalRoleRowName_1"
$ws.Range("R1").Value = "input_rowSelectionCheckbox_internalRoleRowName_2"
$ws.Range("S1").Value = "link_testSuiteActions_internalRoleLinkName"
$ws.Range("T1").Value = "link_testSuiteActions_project_id"
$ws.Range("U1").Value = "link_testSuiteActions_team_id"
$ws.Range("V1").Value = "link_testSuiteActions_test_suites_id"
$ws.Range("W1").Value = "link_testSuiteDetails_executions_id"
$ws.Range("X1").Value = "link_testSuiteDetails_project_id"
$ws.Range("Y1").Value = "link_testSuiteDetails_team_id"

# -----------------------------------------------------------------------
# 3) Row 2 (S2:Y2) data reshuffle.
#    New order: S<-old T, T<-old U, U<-old W, V<-old Y, W<-old S, X<-old U, Y<-old W
#    i.e. final values:
#      S2 = Test Suite Va...Daños Menores
#      T2 = 1408913
#      U2 = 1389363
#      V2 = 25934562-test-suite-validacion-vehiculos-con-danos-menores
#      W2 = 7
#      X2 = 1408913
#      Y2 = 1389363
#    T2/U2/W2/X2/Y2 are purely-numeric-looking text, so force them through
#    as text (NumberFormat "@" while assigning, then reset the style back
#    to Normal so no stray number-format/style is left behind) to avoid
#    Excel auto-converting them into numeric cells.
# -----------------------------------------------------------------------
$ws.Range("S2").Value = "Test Suite Va...Daños Menores"

$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "1408913"
$ws.Range("T2").Style = "Normal"

$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "1389363"
$ws.Range("U2").Style = "Normal"

$ws.Range("V2").Value = "25934562-test-suite-validacion-vehiculos-con-danos-menores"

$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = "7"
$ws.Range("W2").Style = "Normal"

$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "1408913"
$ws.Range("X2").Style = "Normal"

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "1389363"
$ws.Range("Y2").Style = "Normal"
